$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.507.10"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.549.07"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.45"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.59"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.50%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.94%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.18"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.87%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.40"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.01%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.940.47"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.05"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +7.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.595.17"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.844"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.542.73"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.79"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.32"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.43%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0952"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.08"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "244.00"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.69%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.31"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.47%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.40"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.93"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.76"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +13.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0803"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.62"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.19"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.67%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.99%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.20"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +10.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.92"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.77%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.32"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0299"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.964.31"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.36%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.793.70"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.83%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.193"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.07%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "80.93"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.04"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.70%  "
